$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LISTE_S1")

# Update column A date-like values: rows 3..63, subtract 20000 (e.g. 20170926 -> 20150926)
for ($r = 3; $r -le 63; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -ne $null) {
        $cell.Value2 = $val - 20000
    }
}

# Update the sheet view: clear the scrolled topLeftCell and change the selection to A1:C63
$ws.Activate()
$ws.Range("A1:C63").Select()
